# Append the new Key/Value navigation rows (101-108) to the "Navigation" sheet,
# then refresh the view so F1 is the top-left cell and F3 is selected
# (matches the _FilterDatabase range / dimension growing from F100 to F108).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Navigation")

# New Key (col D) / Value (col E) pairs to append starting at row 101.
$rows = @(
    @("AddContactInfo",      "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[2]/div[1]/div[2]/div/div/div/div/a/div/span"),
    @("AddContactEmail",     "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[2]/div[2]/div[1]/div/div[2]/div/div/input"),
    @("AddAddressInfo",      "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[3]/div[1]/div[2]/div/div/div/div/a/div/span"),
    @("AddressTypeDropdown", "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[3]/div[2]/div[1]/div[1]/div[1]/div/div/select"),
    @("AddressPostalCode",   "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[3]/div[2]/div[1]/div[2]/div[1]/div/div/input"),
    @("AddressAdress1",      "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[3]/div[2]/div[1]/div[3]/div/div/div/input"),
    @("AddressSaveButton",   "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[1]/div[3]/div[2]/div[2]/button[1]"),
    @("SelectASClaimant",    "/html/body/div[2]/div/form/div[2]/div[2]/div[2]/div[2]/div[3]/div/div/div[2]/div/div[2]/button[1]")
)

$startRow = 101
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $rows[$i][0]
    $ws.Cells.Item($r, 5).Value = $rows[$i][1]
}

# Update the view: scroll so column F is (closer to) the top-left, and select F3
# (the commit moved topLeftCell from A90 to F1 and the selection from A101 to F3).
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F3").Select()

# The hidden AutoFilter-range defined name needs to grow along with the data.
$filterName = $wb.Names.Item('Navigation!_FilterDatabase')
$filterName.RefersTo = '=Navigation!$A$1:$F$108'
